# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the four "statut" emoji markers used throughout column A with
# plain-text / simple-symbol equivalents:
#   📘 -> ⚠️
#   📕 -> -3
#   📗 -> ✅
#   📙 -> +3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📗" = "✅"
    "📙" = "+3"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

for ($r = $firstRow; $r -lt ($firstRow + $rowCount); $r++) {
    for ($c = $firstCol; $c -lt ($firstCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($map.ContainsKey($val)) {
            $newVal = $map[$val]
            if ($newVal -eq "-3" -or $newVal -eq "+3") {
                # These look like numbers to Excel, so force the cell to
                # Text before assigning, then clear the formatting again
                # so we don't leave a stray number format behind.
                $cell.NumberFormat = "@"
                $cell.Value2 = $newVal
                $cell.ClearFormats()
            } else {
                $cell.Value2 = $newVal
            }
        }
    }
}
